# Nagyobb import, összevont órák(még csak 2 osztály), SolverConfig
#
# Adds two new columns to the "Órák" sheet:
#   E: Osztály2(összevont óra esetén)  -- second class for merged/combined lessons
#   F: Tömbösítés Azonosító            -- grouping id for lessons that are taught
#                                          together ("blocked")
# and fills in the data that goes with the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Órák")

# New header cells (row 1)
$ws.Range("E1").Value = "Osztály2(összevont óra esetén)"
$ws.Range("F1").Value = "Tömbösítés Azonosító"

# Physics is taught together for 9th and 10th grade -> cross-reference the
# other class on both rows.
$ws.Range("E4").Value = "10th grade"
$ws.Range("E15").Value = "9th grade"

# Math lessons taught in two parallel (bontott) groups get a block id so the
# timetable generator knows rows 2/3 belong together, and rows 12/13 belong
# together.
$ws.Range("F2").Value = "a"
$ws.Range("F3").Value = "a"
$ws.Range("F12").Value = "b"
$ws.Range("F13").Value = "b"

# Keep the selection where the author left it when saving.
$ws.Range("G8").Select()
